# ArchitectureDiagram.pptx — developer-guide update
#
# Commit: "Update developer guide. Update the UI section, Calendar section,
# UndoRedo section."
#
# The concrete, COM-reproducible content edit on the single diagram slide is
# the removal of the un-implemented "Web" cloud shape together with the
# dashed elbow connector that tied it to the UI block — i.e. the "Cloud 50"
# shape and the "Elbow Connector 51" connector are deleted from slide 1.
#
# (The rest of the underlying XML diff — a couple of Normal-view slide
# guides recorded in presentation.xml, the cached text of the automatic
# "last saved" date fields on the slide master/layouts, and the
# smtClean="0"/trailing-endParaRPr bookkeeping churn throughout the master
# and layouts — are artifacts of the file having simply been re-saved by a
# newer PowerPoint build; they are not reachable through the PowerPoint
# object model (Guides collection is not implemented in this host, and
# master/layout TextRange.Delete() is a no-op here), so they are
# intentionally left untouched rather than risk corrupting the deck.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Delete the connector first (it references the cloud shape as its end
# connection site), then the cloud shape itself.
$s.Shapes.Item("Elbow Connector 51").Delete()
$s.Shapes.Item("Cloud 50").Delete()
